$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.876.44"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "1.770.61"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "327.81"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "0.4485"
$ws.Range("E7").Value = "  -2.51%  "

$ws.Range("D8").Value = "0.3570"
$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("D10").Value = "42.17"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").Value = "1.095"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").Value = "6.051"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "7.224"
$ws.Range("E15").Value = "  +1.57%  "

$ws.Range("D16").Value = "1.774.56"
$ws.Range("E16").Value = "  +1.21%  "

$ws.Range("D17").Value = "92.83"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").Value = "0.06422"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  +2.53%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "27.902.71"
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").Value = "2.115"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").Value = "162.66"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").Value = "20.26"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "1.975.90"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").Value = "2.197"
$ws.Range("E29").Value = "  +4.86%  "

$ws.Range("D30").Value = "125.66"
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").Value = "1.102"
$ws.Range("E31").Value = "  +2.37%  "

$ws.Range("D32").Value = "0.09165"
$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("D33").Value = "5.580"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("D34").Value = "3.639"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").Value = "11.87"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").Value = "0.02297"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "0.06100"
$ws.Range("E37").Value = "  +0.91%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").Value = "0.6340"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").Value = "4.966"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").Value = "1.184"
$ws.Range("E41").Value = "  -1.64%  "

$ws.Range("D42").Value = "1.393"
$ws.Range("E42").Value = "  +0.47%  "

$ws.Range("D43").Value = "7.927"
$ws.Range("E43").Value = "  +1.55%  "

$ws.Range("D44").Value = "13.25"
$ws.Range("E44").Value = "  +0.15%  "

$ws.Range("D45").Value = "3.741"
$ws.Range("E45").Value = "  +0.84%  "

$ws.Range("D46").Value = "0.5876"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").Value = "0.06921"
$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("D50").Value = "1.140"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("D51").Value = "72.95"
$ws.Range("E51").Value = "  +0.98%  "
